$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column for rows 2-15 from 45208 to 45212
$ws.Range("C2:C15").Value = 45212
